$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 06:52"

# India (row 25) - updated case counts
$ws.Range("B25").Value = 9240
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 1096
$ws.Range("E25").Value = 7813

# Pakistan (row 36) - updated case counts
$ws.Range("B36").Value = 5374
$ws.Range("C36").Value = 144
$ws.Range("D36").Value = 1095
$ws.Range("E36").Value = 4186
$ws.Range("F36").Value = 44
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 93

# Rows 105-107 reorder: Bolivia now appears before Mauricio/Nigeria, each
# keeping/receiving updated case counts (Bolivia gets fresh data, Mauricio
# and Nigeria inherit the previous rows' figures as they shift down).
$ws.Range("A105").Value = "Bolivia"
$ws.Range("B105").Value = 330
$ws.Range("C105").Value = 30
$ws.Range("D105").Value = 2
$ws.Range("E105").Value = 301
$ws.Range("F105").Value = 3
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = 27

$ws.Range("A106").Value = "Mauricio"
$ws.Range("B106").Value = 324
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 42
$ws.Range("E106").Value = 273
$ws.Range("F106").Value = 3
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 9

$ws.Range("A107").Value = "Nigeria"
$ws.Range("B107").Value = 323
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 85
$ws.Range("E107").Value = 228
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 10
